# Daily attendance processing - 2026-01-08 21:34:53
# For every row in the "Recorded By" column (G), if the value begins with
# "System, " (a leading System entry in the recorder list), move that
# "System" token from the front of the list to the back, preserving the
# order/casing of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "System, "
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ($text -ne $null -and $text.StartsWith($prefix)) {
        $rest = $text.Substring($prefix.Length)
        $newValue = $rest + ", System"
        $cell.Value = $newValue
    }
}
